$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range extents
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Swap the contents of column C (codeforiati:group-name) and column D (codeforiati:group-code)
# so that column C becomes the code and column D becomes the name, for every row
# including the header row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value2
    $dValue = $dCell.Value2

    $cCell.Value2 = $dValue
    $dCell.Value2 = $cValue
}
